# Update dSF (column F) values to reflect a repull/recalculation of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = -6
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = -3
